# "Generate Report for Handback"
#
# This localization-status report gets a handback pass recorded:
#  - the status text changes from "Ready for handoff" to
#    "Handed back: in sync with en-US" on the Overview sheet
#  - the per-language sheets (zh-cn / de-de) gain two new columns
#    (F = Latest Target File, G = Latest Handback File) with links to the
#    same files already referenced by the Source File (A) / Handoff File (D)
#    hyperlinks
#  - the Latest Handback DateTime column (H) is populated with the actual
#    handback timestamps instead of the "never happened" placeholder

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Blue color used by the existing hyperlink-styled cells (FF6495ED), and the
# "single" underline style, so the new cells visually match the existing
# handoff/source-file hyperlink cells on each row.
$hyperlinkColor = 15570276
$hyperlinkUnderline = 2

# 1. Overview sheet: status moves from "Ready for handoff" to "Handed back"
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

function Add-HandbackColumns {
    param($ws, $mdTarget, $xlfTarget, $xlfDisplay, $handbackDateTime)

    $mdDisplay = "6ea0face-9cb0-45d6-8b7d-dfac503be676.md"

    foreach ($row in 2,3) {
        $fCell = $ws.Range("F$row")
        $gCell = $ws.Range("G$row")

        $fCell.Value = $mdDisplay
        $ws.Hyperlinks.Add($fCell, $mdTarget, "", "", $mdDisplay) | Out-Null
        $fCell.Font.Underline = $hyperlinkUnderline
        $fCell.Font.Color = $hyperlinkColor

        $gCell.Value = $xlfDisplay
        $ws.Hyperlinks.Add($gCell, $xlfTarget, "", "", $xlfDisplay) | Out-Null
        $gCell.Font.Underline = $hyperlinkUnderline
        $gCell.Font.Color = $hyperlinkColor
    }

    # Latest Handback DateTime now has a real timestamp instead of
    # "0001-01-01 00:00:00"
    $ws.Range("H2").Value = $handbackDateTime
    $ws.Range("H3").Value = $handbackDateTime
}

# 2. zh-cn sheet: new Latest Target File / Latest Handback File columns and
#    handback timestamp
Add-HandbackColumns $zhcn `
    "https://github.com/OpenLocalizationTest/oltest/blob/4d3d64465abe5cbd067f44a875722fc9c6d0e30a/e2e/6ea0face-9cb0-45d6-8b7d-dfac503be676.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b6e22d7e2b43889d2ec241bd4a09035f32ec12f6/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/6ea0face-9cb0-45d6-8b7d-dfac503be676.d072f7e8a3d914fadbfd5d199cad32120bbf389b.zh-cn.xlf" `
    "6ea0face-9cb0-45d6-8b7d-dfac503be676.d072f7e8a3d914fadbfd5d199cad32120bbf389b.zh-cn.xlf" `
    "2016-03-18 17:14:14"

# 3. de-de sheet: new Latest Target File / Latest Handback File columns and
#    handback timestamp
Add-HandbackColumns $dede `
    "https://github.com/OpenLocalizationTest/oltest/blob/4d3d64465abe5cbd067f44a875722fc9c6d0e30a/e2e/6ea0face-9cb0-45d6-8b7d-dfac503be676.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4941d98b685512111fd2ccb3888451b7436c9902/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/6ea0face-9cb0-45d6-8b7d-dfac503be676.d072f7e8a3d914fadbfd5d199cad32120bbf389b.de-de.xlf" `
    "6ea0face-9cb0-45d6-8b7d-dfac503be676.d072f7e8a3d914fadbfd5d199cad32120bbf389b.de-de.xlf" `
    "2016-03-18 17:14:27"

Write-Host "Handback report generated."
